# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row => new F value
$exhibitionUpdates = @{
    6  = 12013
    8  = 79
    9  = 11738
    10 = 4741
    11 = 527
    12 = 69
    14 = 84
    15 = 921
    18 = 54
    19 = 5216
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row => new F value
$allTypesUpdates = @{
    8  = 12013
    10 = 79
    11 = 11738
    12 = 4741
    13 = 527
    14 = 69
    16 = 84
    17 = 921
    20 = 54
    21 = 5216
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
